# Team Roles status update:
# - Center the Status column (F) values
# - Rows that were "-" and are now sworn in get the "Sworn in" label in red
# - Rows still "In Progress" get highlighted in blue

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center-align the whole Status data range first.
$ws.Range("F2:F37").HorizontalAlignment = -4108

# Rows still underway -> keep "In Progress" but call it out in blue.
$inProgressRows = @(32, 35)
foreach ($r in $inProgressRows) {
    $cell = $ws.Range("F$r")
    $cell.Font.Color = 12611584
}

# Rows whose term has now started -> mark "Sworn in" in red.
$swornInRows = @(3, 9, 16, 21, 25)
foreach ($r in $swornInRows) {
    $cell = $ws.Range("F$r")
    $cell.Value = "Sworn in"
    $cell.Font.Color = 192
}

# Restore the cursor position shown in the saved workbook.
$ws.Range("H9").Select()
